$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Annotation scores (Clear, Assertive, Cautious, Optimistic, Specific, Relevant)
# for rows 2-16, columns E-J.
$values = @(
    @(2,2,2,2,2,2),  # row 2
    @(2,2,2,2,2,2),  # row 3
    @(2,2,2,2,2,2),  # row 4
    @(2,2,2,2,2,2),  # row 5
    @(2,2,1,1,2,2),  # row 6
    @(2,2,1,1,2,2),  # row 7
    @(2,2,1,2,1,2),  # row 8
    @(2,2,1,2,2,2),  # row 9
    @(2,2,1,1,2,2),  # row 10
    @(2,2,2,2,2,2),  # row 11
    @(2,2,2,1,1,2),  # row 12
    @(2,2,1,2,2,2),  # row 13
    @(1,2,1,1,2,2),  # row 14
    @(2,1,1,2,2,2),  # row 15
    @(2,2,1,1,2,2)   # row 16
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = 5 + $j   # column E = 5
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Scroll/freeze pane adjustments matching the author's final view state.
$win = $ws.Application.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("C1").Select()
$win.ScrollRow = 1
$win.ScrollColumn = 3
$ws.Range("F17").Select()
